$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Manual")

$ws.Range("H3").Value = ">3% lost (by weight) @cite gellert2015signs"
$ws.Range("L3").Value = ">5% (by weight) lost @cite gellert2015signs"
$ws.Range("P3").Value = ">10% lost (by weight) @cite gellert2015signs"

$ws.Range("J29").Select()
